$wb = $excel.ActiveWorkbook

# --- Rename sheet 2 -------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "GET_last_login"

# --- Sheet2 "GET_last_login" cell edits -----------------------------------
# Row 2: new test case data (was "Get Impact by module guide." / "Get Module Guide" / impacts uri)
$ws2.Range("C2").Value = "Get Last Login"
$ws2.Range("D2").Value = "Get valid userId for Last Login"
$ws2.Range("F2").Value = "/activity/v1/users/{userId}/last-logins"

# Rows 3-5: clear out the leftover GET/status/response cells, keep TCID + Run(N)
$ws2.Range("E3").ClearContents()
$ws2.Range("F3").ClearContents()
$ws2.Range("H3").ClearContents()
$ws2.Range("J3").ClearContents()

$ws2.Range("E4").ClearContents()
$ws2.Range("F4").ClearContents()
$ws2.Range("H4").ClearContents()
$ws2.Range("J4").ClearContents()

$ws2.Range("E5").ClearContents()
$ws2.Range("F5").ClearContents()
$ws2.Range("H5").ClearContents()
$ws2.Range("J5").ClearContents()

# --- Selections / active sheet --------------------------------------------
# Sheet1 is no longer the active tab; update its lingering selection.
[void]$ws1.Activate()
[void]$ws1.Range("E6").Select()

# Sheet2 becomes the active/selected tab.
[void]$ws2.Activate()
[void]$ws2.Range("D8").Select()
